$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(4, "Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .", "Aircraft ActiveTrack available at max speed", "0-5", "Missing")
    ,@(4, "Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .", "When exceeding nnn, Obstacle Avoidance is not available", "7-14", "Missing")
    ,@(4, "Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .", "Aircraft ActiveTrack available at max", "0-4", "False")
    ,@(4, "Aircraft ActiveTrack available at max speed . When exceeding nnn, Obstacle Avoidance is not available .", "When exceeding nnn,", "7-9", "False")
    ,@(13, "Aircraft in Attitude mode . Unable to hover . Fly with caution .", "Aircraft in Attitude mode", "0-3", "Missing")
    ,@(13, "Aircraft in Attitude mode . Unable to hover . Fly with caution .", "Unable to hover", "5-7", "Missing")
    ,@(13, "Aircraft in Attitude mode . Unable to hover . Fly with caution .", "Aircraft in Attitude", "0-2", "False")
    ,@(19, "Aircraft is tilted , please keep the aircraft stationary and level before flight .", "Aircraft is tilted", "0-2", "Missing")
    ,@(27, "Altitude Less than nnn . Cannot enable Follow Me Mode .", "Cannot enable Follow Me Mode", "5-9", "Missing")
    ,@(29, "Another aircraft is dangerously close , please descend to a safer altitude .", "Another aircraft is dangerously close", "0-4", "Missing")
    ,@(48, "Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance .", "If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist", "7-23", "Missing")
    ,@(48, "Check whether propellers are installed correctly . If the propellers are installed correctly and the aircraft still cannot takeoff, a motor error may exist . Contact DJI Support for assistance .", "If the propellers are installed correctly and the aircraft still cannot takeoff,", "7-18", "False")
    ,@(50, "Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff .", "2. Calibrate Compass Before Takeoff", "35-39", "Missing")
    ,@(50, "Compass abnormal . Solution: 1. Ensure there are no magnets or metal objects near the aircraft . The ground or walls may contain metal . Move away from sources of interference before attempting flight . 2. Calibrate Compass Before Takeoff .", "Calibrate Compass Before Takeoff", "36-39", "False")
    ,@(70, "Downward ambient light too low . Obstacle avoidance unavailable . Fly with caution . Backward ambient light too low . Backward obstacle avoidance unavailable . Only infrared sensors available . Fly with caution .", "Backward ambient light too low", "14-18", "Missing")
    ,@(77, "Exiting GPS mode : Unknown Error .", "Unknown Error", "4-5", "Missing")
    ,@(91, "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn .", "GEO Zone Info: The target area is in an Altitude Zone", "0-10", "Missing")
    ,@(91, "GEO Zone Info: The target area is in an Altitude Zone . Flight altitude restricted to nnn .", "GEO Zone Info:", "0-2", "False")
    ,@(92, "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant) . Fly with caution .", "GEO: You are in a Warning Zone (Airport Class Airspace Unpaved Airports Power Plant)", "0-13", "Missing")
    ,@(95, "Gimbal Recenter .", "Gimbal Recenter", "0-1", "Missing")
    ,@(123, "Motor Obstructed . Propulsion output is limited to ensure the health of the battery .", "Propulsion output is limited to ensure the health of the battery", "3-13", "Missing")
    ,@(126, "No GPS . Return to Home failed . Return the aircraft to home manually .", "Return to Home failed", "3-6", "Missing")
    ,@(131, "Panorama Captured Successfully .", "Panorama Captured Successfully", "0-2", "Missing")
    ,@(131, "Panorama Captured Successfully .", "Panorama Captured", "0-1", "False")
    ,@(147, "SD card write speed is too slow . Not suitable for shooting a 4K video .", "Not suitable for shooting a 4K video", "8-14", "Missing")
    ,@(157, "Take off and retry this function .", "Take off and retry this function", "0-5", "Missing")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $errType = $rec[4]
    # "True"/"False" would auto-coerce to native Booleans; prefix with an
    # apostrophe so they stay text, matching the source inlineStr cells.
    $apos = [char]39
    if ($errType -eq "True" -or $errType -eq "False") {
        $ws.Cells.Item($row, 5).Value = "$apos$errType"
    } else {
        $ws.Cells.Item($row, 5).Value = $errType
    }
}

